# Project_Grading_Rubric_Checklist.xlsx edit
# "added randomizer to robot picks to have high school level math"
#
# Concrete content changes applied:
#  - D6: 1 -> 2  (robot picks score bumped)
#  - D9, D31, D32: scores cleared (now blank, grader needs to double check)
#  - D17: highlighted red and a "double check" note added in E17
#  - D26 / D27: yellow "double check" highlight removed now that they're resolved
#  - D33 total recalculates automatically via its SUM formula
#  - selection/scroll position left where the user was last working (D9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Score value changes -----------------------------------------------
$ws.Range("D6").Value = 2

# These three cells had their scores removed (still need review), but they
# keep their existing cell formatting (yellow highlight).
$ws.Range("D9").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("D32").ClearContents()

# --- Flag row 17 for a double check -------------------------------------
$ws.Range("D17").Interior.Color = 255
$ws.Range("E17").Value = "double check"

# --- Remove the now-resolved "double check" highlight on rows 26 & 27 ---
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"

# --- Leave the view where the user left off editing ---------------------
$ws.Range("D9").Select()
